$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price/volume text is preserved exactly as text
# (Excel would otherwise auto-convert strings like "239.00" or "1.0000" to numbers,
# dropping trailing zeros). Apply Text format to the whole data range first.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '29.367.27'
$ws.Range('E2').Value = '  +0.37%  '
$ws.Range('D3').Value = '1.870.53'
$ws.Range('E3').Value = '  +0.56%  '
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '0.7075'
$ws.Range('E5').Value = '  +0.45%  '
$ws.Range('D6').Value = '239.00'
$ws.Range('E6').Value = '  +0.52%  '
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = '0.07782'
$ws.Range('E8').Value = '  -5.54%  '
$ws.Range('D9').Value = '0.3070'
$ws.Range('E9').Value = '  +0.97%  '
$ws.Range('D10').Value = '25.17'
$ws.Range('E10').Value = '  +8.04%  '
$ws.Range('D11').Value = '0.08201'
$ws.Range('E11').Value = '  +0.09%  '
$ws.Range('D12').Value = '1.876.86'
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('D13').Value = '5.248'
$ws.Range('E13').Value = '  +1.34%  '
$ws.Range('D14').Value = '0.7224'
$ws.Range('E14').Value = '  +0.87%  '
$ws.Range('D15').Value = '89.42'
$ws.Range('E15').Value = '  +0.26%  '
$ws.Range('D16').Value = '29.457.41'
$ws.Range('E16').Value = '  +0.59%  '
$ws.Range('D17').Value = '5.816'
$ws.Range('E17').Value = '  +0.68%  '
$ws.Range('D18').Value = '242.02'
$ws.Range('E18').Value = '  +2.02%  '
$ws.Range('D19').Value = '0.000007827'
$ws.Range('E19').Value = '  -0.28%  '
$ws.Range('D20').Value = '13.29'
$ws.Range('E20').Value = '  -0.71%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.139.42'
$ws.Range('E21').Value = '  +0.94%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('D23').Value = '1.0000'
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').Value = '7.746'
$ws.Range('E24').Value = '  +3.80%  '
$ws.Range('D25').Value = '162.30'
$ws.Range('E25').Value = '  +0.23%  '
$ws.Range('D26').Value = '8.966'
$ws.Range('E26').Value = '  -0.24%  '
$ws.Range('D27').Value = '0.1458'
$ws.Range('E27').Value = '  +0.97%  '
$ws.Range('D28').Value = '18.19'
$ws.Range('E28').Value = '  +0.41%  '
$ws.Range('D29').Value = '1.926'
$ws.Range('E29').Value = '  -2.06%  '
$ws.Range('D30').Value = '1.366'
$ws.Range('E30').Value = '  -5.12%  '
$ws.Range('D31').Value = '1.518'
$ws.Range('E31').Value = '  +2.35%  '
$ws.Range('D32').Value = '4.320'
$ws.Range('E32').Value = '  -2.38%  '
$ws.Range('D33').Value = '4.059'
$ws.Range('E33').Value = '  -0.07%  '
$ws.Range('D34').Value = '0.05218'
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('D35').Value = '1.192'
$ws.Range('E35').Value = '  +1.80%  '
$ws.Range('D36').Value = '0.7190'
$ws.Range('E36').Value = '  +1.49%  '
$ws.Range('D37').Value = '1.006'
$ws.Range('E37').Value = '  +0.48%  '
$ws.Range('D38').Value = '2.676'
$ws.Range('E38').Value = '  +0.32%  '
$ws.Range('D39').Value = '0.01856'
$ws.Range('E39').Value = '  +0.33%  '
$ws.Range('D40').Value = '2.703'
$ws.Range('E40').Value = '  -0.81%  '
$ws.Range('D41').Value = '1.176.80'
$ws.Range('E41').Value = '  +3.45%  '
$ws.Range('D42').Value = '0.9171'
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').Value = '5.992'
$ws.Range('E43').Value = '  +0.44%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').Value = '0.4295'
$ws.Range('E44').Value = '  +0.19%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = '71.37'
$ws.Range('E45').Value = '  +0.68%  '
$ws.Range('E46').Value = '  +0.22%  '
$ws.Range('D47').Value = '102.32'
$ws.Range('E47').Value = '  -0.33%  '
$ws.Range('D48').Value = '0.5309'
$ws.Range('E48').Value = '  -2.09%  '
$ws.Range('D49').Value = '1.760'
$ws.Range('E49').Value = '  -0.80%  '
$ws.Range('D50').Value = '9.204'
$ws.Range('E50').Value = '  +0.24%  '
$ws.Range('D51').Value = '7.060'
$ws.Range('E51').Value = '  +1.23%  '
